# Auto-generated Excel COM-interop script
# Applies cached-value corrections to the Adamantoise_Profits workbook
# (static market-price snapshot values refreshed by the scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1657.5714
$ws.Range("I19").Value = 1871.8
$ws.Range("K19").Value = 1871.8
$ws.Range("M19").Value = -1696.8
$ws.Range("H43").Value = 7166.3335
$ws.Range("I43").Value = 4999.5
$ws.Range("K43").Value = 4999.5
$ws.Range("M43").Value = -4930.5
$ws.Range("H44").Value = 38283.332
$ws.Range("J44").Value = 38283.332
$ws.Range("L44").Value = 38283.332
$ws.Range("N44").Value = -39207.332
$ws.Range("H62").Value = 7198.8
$ws.Range("I62").Value = 7248.5
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 7248.5
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -6624.5
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 7198.8
$ws.Range("I65").Value = 7248.5
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 36242.5
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -33122.5
$ws.Range("N65").Value = -41240
$ws.Range("H70").Value = 83344340
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 83344340
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 250033020
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -250033560
$ws.Range("H73").Value = 83344340
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 83344340
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 250033020
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -250034892
$ws.Range("H86").Value = 150074860
$ws.Range("I86").Value = 200011790
$ws.Range("K86").Value = 200011790
$ws.Range("M86").Value = -200010667
$ws.Range("H88").Value = 4038.2307
$ws.Range("J88").Value = 4122.3335
$ws.Range("L88").Value = 4122.3335
$ws.Range("N88").Value = -4934.3335
$ws.Range("H89").Value = 150074860
$ws.Range("I89").Value = 200011790
$ws.Range("K89").Value = 1000058950
$ws.Range("M89").Value = -1000053334
$ws.Range("H91").Value = 4038.2307
$ws.Range("J91").Value = 4122.3335
$ws.Range("L91").Value = 4122.3335
$ws.Range("N91").Value = -6930.3335
$ws.Range("H98").Value = 884.05
$ws.Range("I98").Value = 482.33334
$ws.Range("K98").Value = 482.33334
$ws.Range("M98").Value = 1015.66666
$ws.Range("H115").Value = 906.2
$ws.Range("I115").Value = 907
$ws.Range("J115").Value = 895
$ws.Range("K115").Value = 2721
$ws.Range("L115").Value = 2685
$ws.Range("M115").Value = -1154
$ws.Range("N115").Value = -5819
$ws.Range("H122").Value = 884.05
$ws.Range("I122").Value = 482.33334
$ws.Range("K122").Value = 1447.00002
$ws.Range("M122").Value = 1002.99998
$ws.Range("H138").Value = 1814.74
$ws.Range("J138").Value = 2426.0984
$ws.Range("L138").Value = 7278.2952
$ws.Range("N138").Value = -17558.2952

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19570116
$ws.Range("I32").Value = 22765572
$ws.Range("K32").Value = 22765572
$ws.Range("M32").Value = -22765285
$ws.Range("H45").Value = 2214.2307
$ws.Range("I45").Value = 1253.8889
$ws.Range("K45").Value = 1253.8889
$ws.Range("M45").Value = -876.8888999999999
$ws.Range("H63").Value = 4028.0952
$ws.Range("I63").Value = 2323.625
$ws.Range("K63").Value = 2323.625
$ws.Range("M63").Value = -1637.625
$ws.Range("H66").Value = 4028.0952
$ws.Range("I66").Value = 2323.625
$ws.Range("K66").Value = 11618.125
$ws.Range("M66").Value = -8186.125
$ws.Range("H74").Value = 2180.913
$ws.Range("I74").Value = 2007.1555
$ws.Range("K74").Value = 2007.1555
$ws.Range("M74").Value = -1133.1555
$ws.Range("H77").Value = 2180.913
$ws.Range("I77").Value = 2007.1555
$ws.Range("K77").Value = 10035.7775
$ws.Range("M77").Value = -5667.7775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2783.111
$ws.Range("I105").Value = 2463.818
$ws.Range("J105").Value = 3284.8572
$ws.Range("K105").Value = 2463.818
$ws.Range("L105").Value = 3284.8572
$ws.Range("M105").Value = -716.8180000000002
$ws.Range("N105").Value = -6778.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 774999.75
$ws.Range("I4").Value = 500000
$ws.Range("K4").Value = 500000
$ws.Range("M4").Value = -499888
$ws.Range("H31").Value = 4604.773
$ws.Range("J31").Value = 7045.394
$ws.Range("L31").Value = 7045.394
$ws.Range("N31").Value = -7635.394
$ws.Range("H34").Value = 4604.773
$ws.Range("J34").Value = 7045.394
$ws.Range("L34").Value = 7045.394
$ws.Range("N34").Value = -7449.394
$ws.Range("H99").Value = 3205.3333
$ws.Range("I99").Value = 3052.2222
$ws.Range("K99").Value = 3052.2222
$ws.Range("M99").Value = -1554.2222
$ws.Range("H122").Value = 2945521.8
$ws.Range("I122").Value = 4351674.5
$ws.Range("J122").Value = 5384.273
$ws.Range("K122").Value = 13055023.5
$ws.Range("L122").Value = 16152.819
$ws.Range("M122").Value = -13052573.5
$ws.Range("N122").Value = -21052.819
$ws.Range("H126").Value = 3205.3333
$ws.Range("I126").Value = 3052.2222
$ws.Range("K126").Value = 9156.6666
$ws.Range("M126").Value = -6686.6666
$ws.Range("H132").Value = 3514.158
$ws.Range("I132").Value = 3294.0688
$ws.Range("K132").Value = 9882.206399999999
$ws.Range("M132").Value = -7352.206399999999
$ws.Range("H134").Value = 3580.625
$ws.Range("I134").Value = 3862.4
$ws.Range("K134").Value = 11587.2
$ws.Range("M134").Value = -9052.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 970.7692
$ws.Range("I132").Value = 1124.5
$ws.Range("J132").Value = 942.8182
$ws.Range("K132").Value = 10120.5
$ws.Range("L132").Value = 8485.363800000001
$ws.Range("M132").Value = -7590.5
$ws.Range("N132").Value = -13545.3638
$ws.Range("H133").Value = 5000
$ws.Range("I133").Value = 5000
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -9940
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 3006
$ws.Range("I21").Value = 3006
$ws.Range("K21").Value = 3006
$ws.Range("M21").Value = -2833
$ws.Range("H24").Value = 100028250
$ws.Range("I24").Value = 333346020
$ws.Range("J24").Value = 34919.715
$ws.Range("K24").Value = 333346020
$ws.Range("L24").Value = 34919.715
$ws.Range("M24").Value = -333345847
$ws.Range("N24").Value = -35265.715
$ws.Range("H30").Value = 3006
$ws.Range("I30").Value = 3006
$ws.Range("K30").Value = 3006
$ws.Range("M30").Value = -2901
$ws.Range("H95").Value = 118000
$ws.Range("J95").Value = 118000
$ws.Range("L95").Value = 118000
$ws.Range("N95").Value = -123492
$ws.Range("H126").Value = 2918.75
$ws.Range("I126").Value = 2759
$ws.Range("K126").Value = 8277
$ws.Range("M126").Value = -5807

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H55").Value = 251.5
$ws.Range("I55").Value = 198.125
$ws.Range("K55").Value = 198.125
$ws.Range("M55").Value = -25.125
$ws.Range("H132").Value = 1005182.25
$ws.Range("I132").Value = 1575572.2
$ws.Range("K132").Value = 4726716.6
$ws.Range("M132").Value = -4724186.6
$ws.Range("H135").Value = 51529.293
$ws.Range("J135").Value = 51529.293
$ws.Range("L135").Value = 51529.293
$ws.Range("N135").Value = -61669.293

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5506.143
$ws.Range("J62").Value = 6408.6
$ws.Range("L62").Value = 6408.6
$ws.Range("N62").Value = -7656.6
$ws.Range("H65").Value = 5506.143
$ws.Range("J65").Value = 6408.6
$ws.Range("L65").Value = 32043
$ws.Range("N65").Value = -38283
$ws.Range("H126").Value = 4290.9473
$ws.Range("J126").Value = 2487.1428
$ws.Range("L126").Value = 7461.428400000001
$ws.Range("N126").Value = -12401.4284
$ws.Range("H132").Value = 22986.271
$ws.Range("I132").Value = 26054.537
$ws.Range("K132").Value = 78163.611
$ws.Range("M132").Value = -75633.611
